# Updated cryptos list — refresh price / 1h-volume columns (and resync two
# rank swaps: Kaspa<->Toncoin at rows 29/30, Stacks<->Monero at rows 42/43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $value) {
    # Cells in this sheet are plain inline/shared strings with the default
    # (unstyled) format, even though many values look numeric
    # ("355.92", "13.47", ...). Assigning such a string via COM lets Excel
    # auto-convert it to a real number, which would change the stored cell
    # type. Prefixing with an apostrophe forces text entry (like typing it
    # in the UI); ClearFormats() then drops the resulting "Text"/quote-prefix
    # style so the cell's style index goes back to the original (unstyled).
    $trimmed = $value.Trim()
    $needsApostrophe = $trimmed -match '^[+-]?\d+(\.\d+)?$'

    $range = $ws.Range($cellRef)
    if ($needsApostrophe) {
        $range.Value = "'" + $value
        $range.ClearFormats()
    } else {
        $range.Value = $value
    }
}

# row 2 - Bitcoin
Set-Text "D2" "52.259.73"
Set-Text "E2" "  -0.11%  "

# row 3 - Ethereum
Set-Text "D3" "2.826.61"
Set-Text "E3" "  +1.03%  "

# row 4 - TetherUSD
Set-Text "E4" "  +0.03%  "

# row 5 - BNB
Set-Text "D5" "355.92"
Set-Text "E5" "  +2.54%  "

# row 6 - Solana
Set-Text "D6" "112.24"
Set-Text "E6" "  -3.72%  "

# row 7 - XRP
Set-Text "E7" "  +3.34%  "

# row 8 - USDC
Set-Text "E8" "  +0.05%  "

# row 9 - Cardano
Set-Text "E9" "  +0.85%  "

# row 10 - Avalanche
Set-Text "D10" "40.97"
Set-Text "E10" "  -4.79%  "

# row 11 - Dogecoin
Set-Text "D11" "0.0864"
Set-Text "E11" "  +0.62%  "

# row 12 - TRON
Set-Text "E12" "  +0.96%  "

# row 13 - Chainlink
Set-Text "D13" "19.91"
Set-Text "E13" "  -1.12%  "

# row 14 - Polkadot
Set-Text "D14" "7.76"
Set-Text "E14" "  -1.58%  "

# row 15 - Wrapped liquid staked Ether 2.0
Set-Text "D15" "3.268.51"
Set-Text "E15" "  +0.98%  "

# row 16 - Wrapped Ether
Set-Text "D16" "2.834.91"
Set-Text "E16" "  +1.65%  "

# row 17 - Polygon
Set-Text "D17" "0.932"
Set-Text "E17" "  +4.02%  "

# row 18 - Wrapped BTC
Set-Text "D18" "52.129.94"
Set-Text "E18" "  -0.10%  "

# row 19 - Uniswap
Set-Text "E19" "  +4.46%  "

# row 20 - ImmutableX
Set-Text "E20" "  -1.08%  "

# row 21 - Internet Computer (DFINITY)
Set-Text "D21" "13.47"
Set-Text "E21" "  +0.09%  "

# row 22 - Shiba Inu
Set-Text "D22" "0.0₃0996"
Set-Text "E22" "  +1.43%  "

# row 23 - Litecoin
Set-Text "D23" "70.65"
Set-Text "E23" "  +0.62%  "

# row 24 - Bitcoin Cash
Set-Text "D24" "271.80"
Set-Text "E24" "  +0.59%  "

# row 25 - PancakeSwap
Set-Text "D25" "2.80"
Set-Text "E25" "  +1.09%  "

# row 26 - Ethereum Classic
Set-Text "D26" "26.96"
Set-Text "E26" "  +0.87%  "

# row 27 - Dai
Set-Text "E27" "  +0.13%  "

# row 28 - Cosmos
Set-Text "D28" "10.35"
Set-Text "E28" "  +1.18%  "

# row 29 - was Toncoin, now Kaspa (swapped with row 30)
Set-Text "B29" "Kaspa"
Set-Text "C29" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-Text "D29" "0.145"
Set-Text "E29" "  +3.55%  "

# row 30 - was Kaspa, now Toncoin (swapped with row 29)
Set-Text "B30" "Toncoin"
Set-Text "C30" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-Text "D30" "2.23"
Set-Text "E30" "  -0.94%  "

# row 31 - VeChain
Set-Text "D31" "0.0487"
Set-Text "E31" "  +12.02%  "

# row 32 - OKB
Set-Text "D32" "52.74"
Set-Text "E32" "  +5.05%  "

# row 33 - Injective Protocol
Set-Text "D33" "34.89"
Set-Text "E33" "  -0.46%  "

# row 34 - Filecoin
Set-Text "D34" "5.94"
Set-Text "E34" "  +3.88%  "

# row 35 - Render Token
Set-Text "D35" "5.58"
Set-Text "E35" "  +11.56%  "

# row 36 - Hedera
Set-Text "E36" "  +3.42%  "

# row 38 - Lido DAO Token
Set-Text "D38" "3.27"
Set-Text "E38" "  +1.36%  "

# row 39 - ARBITRUM
Set-Text "E39" "  -3.65%  "

# row 40 - Celestia
Set-Text "E40" "  -2.62%  "

# row 41 - Stellar
Set-Text "E41" "  +1.68%  "

# row 42 - was Monero, now Stacks (swapped with row 43)
Set-Text "B42" "Stacks"
Set-Text "C42" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-Text "D42" "2.56"
Set-Text "E42" "  -4.33%  "

# row 43 - was Stacks, now Monero (swapped with row 42)
Set-Text "B43" "Monero"
Set-Text "C43" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-Text "D43" "127.63"
Set-Text "E43" "  -0.27%  "

# row 44 - EnergySwap
Set-Text "D44" "23.18"
Set-Text "E44" "  -1.87%  "

# row 46 - NEAR Protocol
Set-Text "E46" "  -0.04%  "

# row 47 - Maker
Set-Text "D47" "2.087.95"
Set-Text "E47" "  +0.77%  "

# row 48 - ApeX Protocol
Set-Text "E48" "  -3.90%  "

# row 49 - THORChain
Set-Text "D49" "5.96"

# row 50 - SEI
Set-Text "D50" "0.974"
Set-Text "E50" "  -0.17%  "

# row 51 - FraxShare
Set-Text "E51" "  +2.62%  "
